$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.129.92'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.036.52'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '''514.61'
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").Value = '''140.91'
$ws.Range("E6").Value = '  +2.65%  '
$ws.Range("D8").Value = '''0.440'
$ws.Range("E8").Value = '  +2.72%  '
$ws.Range("D9").Value = '''7.20'
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("E11").Value = '  +5.20%  '
$ws.Range("D12").Value = '3.568.98'
$ws.Range("E12").Value = '  +1.95%  '
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("D14").Value = '''26.78'
$ws.Range("E14").Value = '  +3.71%  '
$ws.Range("D15").Value = '''0.0000165'
$ws.Range("E15").Value = '  +4.73%  '
$ws.Range("D16").Value = '57.049.79'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '3.043.76'
$ws.Range("E17").Value = '  +2.16%  '
$ws.Range("D18").Value = '''6.10'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").Value = '''13.32'
$ws.Range("E19").Value = '  +5.54%  '
$ws.Range("D20").Value = '''8.04'
$ws.Range("E20").Value = '  +3.39%  '
$ws.Range("D21").Value = '''332.58'
$ws.Range("E21").Value = '  +3.62%  '
$ws.Range("D22").Value = '''0.998'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '''0.504'
$ws.Range("E23").Value = '  +3.59%  '
$ws.Range("D24").Value = '''65.36'
$ws.Range("E24").Value = '  +2.81%  '
$ws.Range("D25").Value = '3.176.69'
$ws.Range("E25").Value = '  +2.14%  '
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("D27").Value = '''0.163'
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").Value = '0.0₃0904'
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").Value = '''6.69'
$ws.Range("E29").Value = '  +1.99%  '
$ws.Range("D30").Value = '''7.15'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("E32").Value = '  +3.66%  '
$ws.Range("D33").Value = '''20.63'
$ws.Range("E33").Value = '  +2.32%  '
$ws.Range("E34").Value = '  +1.28%  '
$ws.Range("D35").Value = '''153.31'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '''5.92'
$ws.Range("E36").Value = '  +2.97%  '
$ws.Range("E37").Value = '  +2.26%  '
$ws.Range("D38").Value = '''25.01'
$ws.Range("E38").Value = '  +4.31%  '
$ws.Range("D39").Value = '''0.0671'
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("D40").Value = '3.070.00'
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("D41").Value = '''36.96'
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("E42").Value = '  +3.44%  '
$ws.Range("D43").Value = '''0.998'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").Value = '''0.661'
$ws.Range("E44").Value = '  +3.20%  '
$ws.Range("D45").Value = '2.210.85'
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("D47").Value = '''0.955'
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").Value = '''5.99'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").Value = '''20.25'
$ws.Range("E49").Value = '  +6.03%  '
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("E51").Value = '  +11.54%  '
